$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.055.24"
$ws.Range("E2").Value = "  +1.04%  "

$ws.Range("D3").Value = "2.052.68"
$ws.Range("E3").Value = "  -2.92%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").Value = "'248.24"
$ws.Range("E5").Value = "  -2.72%  "

$ws.Range("D6").Value = "'0.655"
$ws.Range("E6").Value = "  -1.77%  "

$ws.Range("E7").Value = "  -0.08%  "

$ws.Range("D8").Value = "'55.16"
$ws.Range("E8").Value = "  +17.14%  "

$ws.Range("D9").Value = "'61.61"
$ws.Range("E9").Value = "  +1.70%  "

$ws.Range("D10").Value = "'0.378"
$ws.Range("E10").Value = "  +0.66%  "

$ws.Range("D11").Value = "'0.0788"
$ws.Range("E11").Value = "  +5.44%  "

$ws.Range("E12").Value = "  +5.42%  "

$ws.Range("D13").Value = "'15.13"
$ws.Range("E13").Value = "  +5.85%  "

$ws.Range("D14").Value = "2.352.47"
$ws.Range("E14").Value = "  -2.98%  "

$ws.Range("E15").Value = "  -2.36%  "

$ws.Range("D16").Value = "'5.21"
$ws.Range("E16").Value = "  +1.67%  "

$ws.Range("D17").Value = "2.054.70"
$ws.Range("E17").Value = "  -2.85%  "

$ws.Range("D18").Value = "36.975.89"
$ws.Range("E18").Value = "  +0.88%  "

$ws.Range("D19").Value = "0.0₃0908"
$ws.Range("E19").Value = "  +8.35%  "

$ws.Range("D20").Value = "'72.31"
$ws.Range("E20").Value = "  -1.97%  "

$ws.Range("D21").Value = "'14.20"
$ws.Range("E21").Value = "  +7.05%  "

$ws.Range("D22").Value = "'5.30"
$ws.Range("E22").Value = "  +1.95%  "

$ws.Range("D23").Value = "'236.63"
$ws.Range("E23").Value = "  -1.83%  "

$ws.Range("E25").Value = "  -3.29%  "

$ws.Range("D26").Value = "'169.19"
$ws.Range("E26").Value = "  -1.85%  "

$ws.Range("D27").Value = "'8.98"
$ws.Range("E27").Value = "  -2.50%  "

$ws.Range("D28").Value = "'20.01"
$ws.Range("E28").Value = "  -7.95%  "

$ws.Range("D29").Value = "'1.96"
$ws.Range("E29").Value = "  -2.93%  "

$ws.Range("E30").Value = "  -0.45%  "

$ws.Range("D31").Value = "'4.52"
$ws.Range("E31").Value = "  +0.50%  "

$ws.Range("D32").Value = "'1.05"
$ws.Range("E32").Value = "  +9.53%  "

$ws.Range("D33").Value = "'0.0621"
$ws.Range("E33").Value = "  +3.37%  "

$ws.Range("D34").Value = "'4.29"
$ws.Range("E34").Value = "  +3.25%  "

$ws.Range("E35").Value = "  -0.09%  "

$ws.Range("D36").Value = "'0.0860"
$ws.Range("E36").Value = "  -9.72%  "

$ws.Range("D37").Value = "'2.25"
$ws.Range("E37").Value = "  -4.38%  "

$ws.Range("D38").Value = "'1.77"
$ws.Range("E38").Value = "  -6.04%  "

$ws.Range("D39").Value = "'1.34"
$ws.Range("E39").Value = "  -0.08%  "

$ws.Range("E40").Value = "  +20.90%  "

$ws.Range("D41").Value = "'18.12"
$ws.Range("E41").Value = "  +12.47%  "

$ws.Range("E42").Value = "  -1.32%  "

$ws.Range("D43").Value = "'1.13"
$ws.Range("E43").Value = "  -4.76%  "

$ws.Range("D44").Value = "'95.65"
$ws.Range("E44").Value = "  -3.28%  "

$ws.Range("B45").Value = "HuobiToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D45").Value = "'2.77"
$ws.Range("E45").Value = "  -0.75%  "

$ws.Range("B46").Value = "FTXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D46").Value = "'4.15"
$ws.Range("E46").Value = "  +40.08%  "

$ws.Range("D47").Value = "'14.09"
$ws.Range("E47").Value = "  -51.35%  "

$ws.Range("D48").Value = "'2.39"
$ws.Range("E48").Value = "  +4.28%  "

$ws.Range("D49").Value = "1.293.98"
$ws.Range("E49").Value = "  -4.49%  "

$ws.Range("D50").Value = "'2.90"
$ws.Range("E50").Value = "  +2.39%  "

$ws.Range("D51").Value = "'6.73"
$ws.Range("E51").Value = "  -6.87%  "
